$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.081.68"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.891.52"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.26"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.56"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.85"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0856"
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.94"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.77"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "3.343.18"
$ws.Range("E15").Value = "  +3.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.01"
$ws.Range("E16").Value = "  +7.31%  "
$ws.Range("D17").Value = "2.931.05"
$ws.Range("E17").Value = "  +5.46%  "
$ws.Range("D18").Value = "52.073.19"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.71"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.35"
$ws.Range("E20").Value = "  +6.08%  "
$ws.Range("E21").Value = "  +8.02%  "
$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.73"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.41"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.48"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.46"
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  +8.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.43"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0944"
$ws.Range("E34").Value = "  +10.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.91"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  +5.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.59"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.65"
$ws.Range("E41").Value = "  +6.52%  "
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.62"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.02"
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.20"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.57"
$ws.Range("E46").Value = "  +4.78%  "
$ws.Range("D47").Value = "2.202.81"
$ws.Range("E47").Value = "  +3.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.51"
$ws.Range("E48").Value = "  +6.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.270"
$ws.Range("E49").Value = "  +21.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.946"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0322"
$ws.Range("E51").Value = "  +10.33%  "
